# Fix: Dynamic Dag dictionary
#
# The "daily" sheet's ticker list had "JPY_USD" — it should be "USD_JPY".
# Fixing that cell is the actual content change; it also shuffles the
# workbook's shared-strings table (the corrected string sorts to the end)
# which shifts the shared-string indices that the other, unchanged ticker
# rows point to — that part is just a side-effect of the rename, handled
# automatically by the engine when we write the new cell value.

$wb = $excel.ActiveWorkbook

$wsDaily = $wb.Worksheets.Item("daily")

# Correct the ticker symbol in A3 from JPY_USD to USD_JPY.
$wsDaily.Range("A3").Value = "USD_JPY"

# Update the selected/active cell on the "daily" sheet.
$wsDaily.Range("E4").Select()

# Best-effort: reflect the saved window geometry (position/size) too.
$excel.Width = 14400
$excel.Height = 7360
$excel.Left = 380
$excel.Top = 380

$wb.Save()
